$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2013888888888889
$ws.Range("C2").Value = 0.5659722222222222
$ws.Range("J2").Value = 0.01388888888888889
$ws.Range("P2").Value = 0.1493055555555556
$ws.Range("S2").Value = 0.06944444444444445
$ws.Range("B3").Value = 0.01694915254237288
$ws.Range("C3").Value = 0.03954802259887006
$ws.Range("J3").Value = 0.01694915254237288
$ws.Range("P3").Value = 0.7062146892655368
$ws.Range("S3").Value = 0.2203389830508475
$ws.Range("J4").Value = 0.02380952380952381
$ws.Range("P4").Value = 0.7380952380952381
$ws.Range("S4").Value = 0.2380952380952381
$ws.Range("P5").Value = 0.6
$ws.Range("S5").Value = 0.4
$ws.Range("B6").Value = 0.04700854700854701
$ws.Range("D6").Value = 0.008547008547008548
$ws.Range("E6").Value = 0.004273504273504274
$ws.Range("F6").Value = 0.0811965811965812
$ws.Range("J6").Value = 0.2478632478632479
$ws.Range("Q6").Value = 0.141025641025641
$ws.Range("R6").Value = 0.05555555555555555
$ws.Range("S6").Value = 0.4145299145299146
$ws.Range("B7").Value = 0.05853658536585366
$ws.Range("D7").Value = 0.02439024390243903
$ws.Range("E7").Value = 0.004878048780487805
$ws.Range("F7").Value = 0.05853658536585366
$ws.Range("J7").Value = 0.1317073170731707
$ws.Range("O7").Value = 0.01951219512195122
$ws.Range("Q7").Value = 0.1170731707317073
$ws.Range("R7").Value = 0.07804878048780488
$ws.Range("S7").Value = 0.5073170731707317
$ws.Range("B8").Value = 0.08921161825726141
$ws.Range("D8").Value = 0.02074688796680498
$ws.Range("E8").Value = 0.002074688796680498
$ws.Range("F8").Value = 0.06639004149377593
$ws.Range("J8").Value = 0.1244813278008299
$ws.Range("O8").Value = 0.02904564315352697
$ws.Range("Q8").Value = 0.2074688796680498
$ws.Range("R8").Value = 0.09336099585062241
$ws.Range("S8").Value = 0.3672199170124482
$ws.Range("B9").Value = 0.1
$ws.Range("D9").Value = 0.007407407407407408
$ws.Range("E9").Value = 0.007407407407407408
$ws.Range("F9").Value = 0.02962962962962963
$ws.Range("J9").Value = 0.1037037037037037
$ws.Range("O9").Value = 0.02592592592592593
$ws.Range("Q9").Value = 0.1666666666666667
$ws.Range("R9").Value = 0.08888888888888889
$ws.Range("S9").Value = 0.4703703703703704
$ws.Range("B10").Value = 0.09368191721132897
$ws.Range("D10").Value = 0.01742919389978214
$ws.Range("F10").Value = 0.06681190994916485
$ws.Range("J10").Value = 0.1336238198983297
$ws.Range("O10").Value = 0.01597676107480029
$ws.Range("Q10").Value = 0.1902687000726216
$ws.Range("R10").Value = 0.0915032679738562
$ws.Range("S10").Value = 0.3907044299201162
$ws.Range("G11").Value = 0.1220238095238095
$ws.Range("J11").Value = 0.07142857142857142
$ws.Range("K11").Value = 0.1845238095238095
$ws.Range("L11").Value = 0.5982142857142857
$ws.Range("S11").Value = 0.02380952380952381
$ws.Range("G12").Value = 0.73
$ws.Range("J12").Value = 0.19
$ws.Range("L12").Value = 0.005
$ws.Range("S12").Value = 0.075
$ws.Range("G13").Value = 0.6857142857142857
$ws.Range("J13").Value = 0.2
$ws.Range("S13").Value = 0.1142857142857143
$ws.Range("F15").Value = 0.02032520325203252
$ws.Range("H15").Value = 0.1747967479674797
$ws.Range("I15").Value = 0.08536585365853659
$ws.Range("J15").Value = 0.3292682926829268
$ws.Range("K15").Value = 0.04878048780487805
$ws.Range("M15").Value = 0.008130081300813009
$ws.Range("O15").Value = 0.06504065040650407
$ws.Range("S15").Value = 0.2682926829268293
$ws.Range("F16").Value = 0.02105263157894737
$ws.Range("H16").Value = 0.1684210526315789
$ws.Range("I16").Value = 0.1210526315789474
$ws.Range("J16").Value = 0.4052631578947368
$ws.Range("K16").Value = 0.09473684210526316
$ws.Range("M16").Value = 0.005263157894736842
$ws.Range("N16").Value = 0.005263157894736842
$ws.Range("O16").Value = 0.04736842105263158
$ws.Range("S16").Value = 0.131578947368421
$ws.Range("F17").Value = 0.02114164904862579
$ws.Range("H17").Value = 0.1627906976744186
$ws.Range("I17").Value = 0.105708245243129
$ws.Range("J17").Value = 0.4207188160676533
$ws.Range("K17").Value = 0.09513742071881606
$ws.Range("M17").Value = 0.01691331923890063
$ws.Range("N17").Value = 0.002114164904862579
$ws.Range("O17").Value = 0.06553911205073996
$ws.Range("S17").Value = 0.1099365750528541
$ws.Range("F18").Value = 0.008888888888888889
$ws.Range("H18").Value = 0.1688888888888889
$ws.Range("I18").Value = 0.1377777777777778
$ws.Range("J18").Value = 0.4488888888888889
$ws.Range("K18").Value = 0.07111111111111111
$ws.Range("M18").Value = 0.008888888888888889
$ws.Range("O18").Value = 0.05333333333333334
$ws.Range("S18").Value = 0.1022222222222222
$ws.Range("F19").Value = 0.01753202966958867
$ws.Range("H19").Value = 0.2016183412002697
$ws.Range("I19").Value = 0.09912339851652056
$ws.Range("J19").Value = 0.3432231962238705
$ws.Range("K19").Value = 0.1200269723533378
$ws.Range("M19").Value = 0.01618341200269724
$ws.Range("N19").Value = 0.0006743088334457181
$ws.Range("O19").Value = 0.06540795684423466
$ws.Range("S19").Value = 0.1362103843560351
